# Popular car model implemented
# Add a new worksheet "PopularCarModels" after the existing "BikeDetails"
# sheet, populate it with a header + 9 popular car model names, and give
# the header cell the same "white on green" look used by the header row
# of the BikeDetails sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet as the LAST sheet (i.e. after BikeDetails).
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "PopularCarModels"

$carModels = @(
    "Popular Car Model",
    "Maruti 800",
    "Maruti Swift Dzire",
    "Maruti Swift",
    "Hyundai I10",
    "Hyundai Santro Xing",
    "Honda City",
    "Toyota Innova",
    "Toyota Fortuner",
    "Mahindra XUV500"
)

for ($i = 0; $i -lt $carModels.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $carModels[$i]
}

# Style the header cell (A1) like the other sheet's header row: white
# text on a solid green fill.
$header = $ws2.Range("A1")
$header.Font.Color = 16777215
$header.Interior.Color = 32768
$header.Interior.Pattern = 1

# Column A should be just wide enough to fit the longest entry.
$ws2.Columns.Item(1).ColumnWidth = 16.338541666666668

# Restore BikeDetails as the active sheet/tab (adding a sheet makes the
# new one active by default).
$ws1.Activate()

Write-Host "PopularCarModels sheet added"
